$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing (pre-edit) data for rows 60 and 61 ---
# These rows will be pushed down to become the new rows 62 and 63
# once we insert two fresh rows below the (updated) rows 60-61.
$orig60 = @()
$orig61 = @()
for ($c = 1; $c -le 20; $c++) {
    $orig60 += ,$ws.Cells.Item(60, $c).Value()
    $orig61 += ,$ws.Cells.Item(61, $c).Value()
}

# --- Insert two new rows at positions 62-63 ---
# This shifts the old row 62 down to row 64, while rows 60-61 stay put.
$ws.Rows("62:63").Insert()

# --- Restore the original (pre-edit) row 60 / row 61 data into the ---
# --- newly created rows 62 and 63                                   ---
for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item(62, $c).Value = $orig60[$c - 1]
    $ws.Cells.Item(63, $c).Value = $orig61[$c - 1]
}
# Preserve the date number format on column D for the restored rows
$ws.Range("D62").NumberFormat = $ws.Range("D60").NumberFormat
$ws.Range("D63").NumberFormat = $ws.Range("D60").NumberFormat

# --- Update row 60: date moves to 2023-03-28 (serial 45013) ---
$ws.Range("D60").Value = 45013

# --- Update row 61: new price entry dated 2023-03-28 ---
$ws.Range("D61").Value = 45013
$ws.Range("M61").Value = 200
$ws.Range("N61").Value = 8000
$ws.Range("O61").Value = 8000
$ws.Range("P61").Value = 8000
$ws.Range("Q61").Value = "`$/caja 18 kilos granel"
$ws.Range("R61").Value = "Provincia de Curicó"
$ws.Range("S61").Value = 444
$ws.Range("T61").Value = 18
